$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 184 (shifts existing rows 184-228 down to 185-229)
$ws.Rows.Item(184).Insert()

# Populate the newly inserted row 184 with a new weekly data point.
# Values mirror the existing "Vega Modelo de Temuco" / Espinaca rows, with a
# new date (2023-01-02 serial 44932) and updated volume/price figures.
$ws.Cells.Item(184, 1).Value = 10
$ws.Cells.Item(184, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(184, 3).Value = "La Araucanía"
$ws.Cells.Item(184, 4).Value = 44932
$ws.Cells.Item(184, 5).Value = 9
$ws.Cells.Item(184, 6).Value = 100112012
$ws.Cells.Item(184, 7).Value = "Espinaca"
$ws.Cells.Item(184, 8).Value = "Sin especificar"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 30
$ws.Cells.Item(184, 11).Value = 10000
$ws.Cells.Item(184, 12).Value = 10000
$ws.Cells.Item(184, 13).Value = 10000
$ws.Cells.Item(184, 14).Value = "$/docena de atados"
$ws.Cells.Item(184, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(184, 16).Value = 3333
$ws.Cells.Item(184, 17).Value = 3
$ws.Cells.Item(184, 18).Value = "Hortaliza"
